# #5: cash & deposit done
# Rebuild the "存款" (deposit) sheet (sheet index 2) so every row follows the
# same normalized column layout used by the other property sheets:
#   bank | deposit_type | currency | owner | total | property_category |
#   category | date | legislator_name | legislator_id | source_file | index
#
# Previously the sheet only had columns B..G (bank/deposit_type/currency/
# owner/<blank>/amount) and row 1 was accidentally a second data row instead
# of a header row. This fixes the header and appends the missing metadata
# columns (H..M) that the other sheets already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Row 1: turn the stray duplicate data row into a real header row ------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"

# New header cells G1:M1 - copy the header style (bold + border) from B1
# first, then fill in the labels.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1:M1").PasteSpecial(-4122) | Out-Null
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

$excel.CutCopyMode = 0

# --- Data rows 2-8 ----------------------------------------------------------
# columns: A index(unchanged) B bank C deposit_type D currency E owner
#          F total G property_category H category I date J legislator_name
#          K legislator_id L source_file M index(dup of A)
$rows = @(
    @{ Row=2;  Index=44; Bank="華南商業銀行台大分行";         Type="活期存款";       Currency="新臺幣"; Total=6480 },
    @{ Row=3;  Index=45; Bank="匯豐（台灣）商業銀行台北分行"; Type="活期存款";       Currency="新臺幣"; Total=672819 },
    @{ Row=4;  Index=46; Bank="台新國際商業銀行板橋分行";     Type="活期存款";       Currency="新臺幣"; Total=8366 },
    @{ Row=5;  Index=47; Bank="遠東國際商業銀行板橋埔墘分行"; Type="活期存款";       Currency="新臺幣"; Total=4535528 },
    @{ Row=6;  Index=48; Bank="台北台大郵局（第23支局）";     Type="中華郵政存簿儲金"; Currency="新臺幣"; Total=1233528 },
    @{ Row=7;  Index=49; Bank="台北台大郵局（第23支局）";     Type="公教優惠儲蓄存款"; Currency="新臺幣"; Total=630441 },
    @{ Row=8;  Index=50; Bank="遠東國際商業銀行板橋埔墘分行"; Type="外幣存款";       Currency="人民幣"; Total=9627.12 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Existing columns B:G already exist in the sheet (same style) -- just
    # overwrite their values / meaning.
    $ws.Cells.Item($row, 2).Value = $r.Bank       # B bank
    $ws.Cells.Item($row, 3).Value = $r.Type       # C deposit_type
    $ws.Cells.Item($row, 4).Value = $r.Currency   # D currency
    $ws.Cells.Item($row, 5).Value = "孫效智"       # E owner (unchanged)
    $ws.Cells.Item($row, 6).Value = $r.Total      # F total (was blank/quantity)

    # New columns G:M - clone the data-row style (border-less, s=2) from an
    # existing data cell (e.g. column B of this row) before assigning values.
    $ws.Range("B$row").Copy() | Out-Null
    $ws.Range("G$row:M$row").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 7).Value  = "deposit"       # G property_category
    $ws.Cells.Item($row, 8).Value  = "normal"        # H category
    $ws.Cells.Item($row, 9).Value  = "2013-12-19"    # I date
    $ws.Cells.Item($row, 10).Value = "楊玉欣"         # J legislator_name
    $ws.Cells.Item($row, 11).Value = 1757            # K legislator_id
    $ws.Cells.Item($row, 12).Value = "tmp7d8c1"      # L source_file
    $ws.Cells.Item($row, 13).Value = $r.Index        # M index
}

$excel.CutCopyMode = 0
